# Auto-generated COM-interop script to apply the "2020-09-27" data update
# to the "Fonds de solidarite, volet 2" regional/classe_effectif sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C, D, E, G can hold numeric-looking text (counts, amounts with
# trailing ".00", zero-padded codes like "00"/"01"/"02"/"03") that must stay
# stored as text, exactly like the source inlineStr cells. Force the text
# number format before assigning so Excel does not silently coerce them to
# numbers (which would drop the padding / trailing zeros).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Phase 1: refresh nombre_aides / montant_total for existing rows 2-87 ---
# (these rows are above the insertion point, so their row numbers do not move)
Set-TextValue $ws.Range("C2") "210"
Set-TextValue $ws.Range("D2") "568228.00"
Set-TextValue $ws.Range("C3") "1109"
Set-TextValue $ws.Range("D3") "3849213.01"
Set-TextValue $ws.Range("C4") "446"
Set-TextValue $ws.Range("D4") "2079230.78"
Set-TextValue $ws.Range("C5") "129"
Set-TextValue $ws.Range("D5") "682128.09"
Set-TextValue $ws.Range("C8") "50"
Set-TextValue $ws.Range("D8") "106000.00"
Set-TextValue $ws.Range("C9") "65"
Set-TextValue $ws.Range("D9") "168607.56"
Set-TextValue $ws.Range("C10") "382"
Set-TextValue $ws.Range("D10") "1408261.71"
Set-TextValue $ws.Range("C11") "158"
Set-TextValue $ws.Range("D11") "752101.77"
Set-TextValue $ws.Range("C12") "41"
Set-TextValue $ws.Range("D12") "260852.00"
Set-TextValue $ws.Range("C14") "17"
Set-TextValue $ws.Range("D14") "43500.00"
Set-TextValue $ws.Range("C15") "108"
Set-TextValue $ws.Range("D15") "304126.38"
Set-TextValue $ws.Range("C16") "446"
Set-TextValue $ws.Range("D16") "1437153.23"
Set-TextValue $ws.Range("C17") "169"
Set-TextValue $ws.Range("D17") "830937.38"
Set-TextValue $ws.Range("C18") "51"
Set-TextValue $ws.Range("D18") "268045.00"
Set-TextValue $ws.Range("C21") "56"
Set-TextValue $ws.Range("D21") "134800.00"
Set-TextValue $ws.Range("C22") "345"
Set-TextValue $ws.Range("D22") "1154026.51"
Set-TextValue $ws.Range("C23") "126"
Set-TextValue $ws.Range("D23") "601660.00"
Set-TextValue $ws.Range("C27") "21"
Set-TextValue $ws.Range("D27") "51500.00"
Set-TextValue $ws.Range("C34") "122"
Set-TextValue $ws.Range("D34") "362673.00"
Set-TextValue $ws.Range("C35") "632"
Set-TextValue $ws.Range("D35") "2255462.10"
Set-TextValue $ws.Range("C36") "276"
Set-TextValue $ws.Range("D36") "1600678.11"
Set-TextValue $ws.Range("C37") "91"
Set-TextValue $ws.Range("D37") "670593.00"
Set-TextValue $ws.Range("C38") "30"
Set-TextValue $ws.Range("D38") "230500.00"
Set-TextValue $ws.Range("C39") "32"
Set-TextValue $ws.Range("D39") "68200.00"
Set-TextValue $ws.Range("C40") "51"
Set-TextValue $ws.Range("D40") "146636.00"
Set-TextValue $ws.Range("C41") "200"
Set-TextValue $ws.Range("D41") "574160.00"
Set-TextValue $ws.Range("C42") "91"
Set-TextValue $ws.Range("D42") "342429.00"
Set-TextValue $ws.Range("C45") "85"
Set-TextValue $ws.Range("D45") "258917.00"
Set-TextValue $ws.Range("C46") "39"
Set-TextValue $ws.Range("D46") "145357.84"
Set-TextValue $ws.Range("C47") "114"
Set-TextValue $ws.Range("D47") "538974.61"
Set-TextValue $ws.Range("C48") "62"
Set-TextValue $ws.Range("D48") "360703.00"
Set-TextValue $ws.Range("C49") "35"
Set-TextValue $ws.Range("D49") "311681.00"
Set-TextValue $ws.Range("C51") "23"
Set-TextValue $ws.Range("D51") "78900.00"
Set-TextValue $ws.Range("C52") "117"
Set-TextValue $ws.Range("D52") "353786.05"
Set-TextValue $ws.Range("C53") "710"
Set-TextValue $ws.Range("D53") "2996539.98"
Set-TextValue $ws.Range("C54") "300"
Set-TextValue $ws.Range("D54") "1521617.74"
Set-TextValue $ws.Range("C55") "116"
Set-TextValue $ws.Range("D55") "807097.18"
Set-TextValue $ws.Range("C58") "859"
Set-TextValue $ws.Range("D58") "2357840.24"
Set-TextValue $ws.Range("C59") "4197"
Set-TextValue $ws.Range("D59") "14280539.02"
Set-TextValue $ws.Range("C60") "2163"
Set-TextValue $ws.Range("D60") "10303623.78"
Set-TextValue $ws.Range("C61") "747"
Set-TextValue $ws.Range("D61") "4084097.06"
Set-TextValue $ws.Range("C62") "160"
Set-TextValue $ws.Range("D62") "1177723.00"
Set-TextValue $ws.Range("C64") "440"
Set-TextValue $ws.Range("D64") "1131579.34"
Set-TextValue $ws.Range("C65") "39"
Set-TextValue $ws.Range("D65") "103961.00"
Set-TextValue $ws.Range("C84") "248"
Set-TextValue $ws.Range("D84") "686982.39"
Set-TextValue $ws.Range("C85") "952"
Set-TextValue $ws.Range("D85") "3390636.16"
Set-TextValue $ws.Range("C86") "363"
Set-TextValue $ws.Range("D86") "1835669.18"
Set-TextValue $ws.Range("C87") "135"
Set-TextValue $ws.Range("D87") "901984.05"

# --- Phase 2: insert the new "Nouvelle-Aquitaine / 20 a 49 salaries" row ---
# Inserting at row 89 pushes the former rows 89-107 down to 90-108 and keeps
# their A/B/E/F/G/H (dispositif/volet/reg/libelle_region/classe_effectif/
# libelle_classe_effectif) values intact automatically.
$ws.Rows(89).Insert()

$ws.Range("A89").Value = "Fonds de solidarité"
$ws.Range("B89").Value = "VOLET2"
Set-TextValue $ws.Range("C89") "3"
Set-TextValue $ws.Range("D89") "65000.00"
Set-TextValue $ws.Range("E89") "75"
$ws.Range("F89").Value = "Nouvelle-Aquitaine"
Set-TextValue $ws.Range("G89") "12"
$ws.Range("H89").Value = "20 à 49 salariés"

# --- Phase 3: refresh nombre_aides / montant_total for the rows that shifted
#     down one position (formerly 89-107, now 90-108) ---
Set-TextValue $ws.Range("C90") "38"
Set-TextValue $ws.Range("D90") "90500.00"
Set-TextValue $ws.Range("C92") "1084"
Set-TextValue $ws.Range("D92") "3540548.59"
Set-TextValue $ws.Range("C98") "434"
Set-TextValue $ws.Range("D98") "1375350.36"
Set-TextValue $ws.Range("C99") "188"
Set-TextValue $ws.Range("D99") "818905.72"
Set-TextValue $ws.Range("C100") "64"
Set-TextValue $ws.Range("D100") "375384.17"
Set-TextValue $ws.Range("C102") "17"
Set-TextValue $ws.Range("D102") "38260.00"
Set-TextValue $ws.Range("C103") "344"
Set-TextValue $ws.Range("D103") "931328.82"
Set-TextValue $ws.Range("C104") "1373"
Set-TextValue $ws.Range("D104") "4515948.55"
Set-TextValue $ws.Range("C105") "511"
Set-TextValue $ws.Range("D105") "2278695.33"
Set-TextValue $ws.Range("C106") "142"
Set-TextValue $ws.Range("D106") "883996.00"
Set-TextValue $ws.Range("C107") "45"
Set-TextValue $ws.Range("D107") "343657.00"
Set-TextValue $ws.Range("C108") "91"
Set-TextValue $ws.Range("D108") "206456.16"
